$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.411.29"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.817.73"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'315.47"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'0.5115"
$ws.Range("E7").Value = "  -4.38%  "
$ws.Range("D8").Value = "'0.3952"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("D9").Value = "'0.08041"
$ws.Range("E9").Value = "  +5.67%  "
$ws.Range("D10").Value = "'41.71"
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("D11").Value = "'1.107"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "'20.97"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "'6.257"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D16").Value = "1.816.48"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "'0.00001139"
$ws.Range("E17").Value = "  +6.17%  "
$ws.Range("D19").Value = "'0.06624"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'6.082"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "28.440.10"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'11.24"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "'2.268"
$ws.Range("E25").Value = "  +3.31%  "
$ws.Range("D26").Value = "'21.11"
$ws.Range("D27").Value = "2.036.50"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "'154.83"
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").Value = "'2.404"
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("D30").Value = "'125.79"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("D31").Value = "'0.1100"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("D33").Value = "'5.761"
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("D34").Value = "'3.651"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "'0.07022"
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("D36").Value = "'0.2225"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "'0.02329"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "'5.199"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").Value = "'8.809"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("D40").Value = "'0.6256"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "'1.396"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "'13.49"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").Value = "'124.86"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").Value = "'1.973"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "'0.06890"
$ws.Range("E51").Value = "  +0.08%  "
